# Remove the notion of pre-installed applications in the build-flow overview
# slide. This is not important at this step; it is explained further in the
# Virtual Device builder.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1. Remove the "Application 1" placeholder box.
$s.Shapes.Item("Rounded Rectangle 10").Delete()

# 2. Remove the "Application N" placeholder box.
$s.Shapes.Item("Rounded Rectangle 71").Delete()

# 3. Remove the down arrow that pointed from the pre-installed application
#    box down into the "VEE Port" box.
$s.Shapes.Item("Down Arrow 30").Delete()

# 4. Remove the dashed rectangle that grouped/highlighted the pre-installed
#    application boxes.
$s.Shapes.Item("Rectangle 21").Delete()

# 5. Remove the "Pre-installed Applications" caption.
$s.Shapes.Item("Rectangle 22").Delete()

# 6. Tidy up the "VEE Port" label: it used to be split across two runs
#    ("EE " and "Port"); merge them into a single run of text.
$veePort = $s.Shapes.Item("Rounded Rectangle 24").TextFrame.TextRange
$veePort.Characters(2, $veePort.Text.Length - 1).Text = "EE Port"

# 7. Tidy up the "Kernel binary" label: it used to be split across two runs
#    ("Kernel " and "binary"); merge them into a single run of text.
$kernelBinary = $s.Shapes.Item("Rounded Rectangle 23").TextFrame.TextRange
$kernelBinary.Characters(1, 13).Text = "Kernel binary"
